$wb = $excel.ActiveWorkbook


# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 2533.8809
$ws.Range("I98").Value = 2321.7932
$ws.Range("K98").Value = 2321.7932
$ws.Range("M98").Value = -823.7932000000001
# Row 122
$ws.Range("H122").Value = 2533.8809
$ws.Range("I122").Value = 2321.7932
$ws.Range("K122").Value = 6965.3796
$ws.Range("M122").Value = -4515.3796
# Row 132
$ws.Range("H132").Value = 3061.2632
$ws.Range("I132").Value = 2854.4666
$ws.Range("J132").Value = 3836.75
$ws.Range("K132").Value = 8563.399800000001
$ws.Range("L132").Value = 11510.25
$ws.Range("M132").Value = -6033.399800000001
$ws.Range("N132").Value = -16570.25
# Row 137
$ws.Range("H137").Value = 4491.1
$ws.Range("I137").Value = 1407.5
$ws.Range("K137").Value = 4222.5
$ws.Range("M137").Value = -1672.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
# Row 32
$ws.Range("H32").Value = 1075.5195
$ws.Range("I32").Value = 698.8143
$ws.Range("K32").Value = 698.8143
$ws.Range("M32").Value = -411.8143
# Row 61
$ws.Range("H61").Value = 6468.95
$ws.Range("I61").Value = 3870.1428
$ws.Range("J61").Value = 12532.833
$ws.Range("K61").Value = 3870.1428
$ws.Range("L61").Value = 12532.833
$ws.Range("M61").Value = -3658.1428
$ws.Range("N61").Value = -12956.833
# Row 97
$ws.Range("H97").Value = 1520.3182
$ws.Range("I97").Value = 698.1177
$ws.Range("J97").Value = 4315.8
$ws.Range("K97").Value = 698.1177
$ws.Range("L97").Value = 4315.8
$ws.Range("M97").Value = -202.1177
$ws.Range("N97").Value = -5307.8
# Row 136
$ws.Range("H136").Value = 6468.95
$ws.Range("I136").Value = 3870.1428
$ws.Range("J136").Value = 12532.833
$ws.Range("K136").Value = 11610.4284
$ws.Range("L136").Value = 37598.499
$ws.Range("M136").Value = -9060.428400000001
$ws.Range("N136").Value = -42698.499

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5380.6177
$ws.Range("I86").Value = 2292.6667
$ws.Range("K86").Value = 2292.6667
$ws.Range("M86").Value = -1169.6667
# Row 89
$ws.Range("H89").Value = 5380.6177
$ws.Range("I89").Value = 2292.6667
$ws.Range("K89").Value = 11463.3335
$ws.Range("M89").Value = -5847.333500000001
# Row 119
$ws.Range("H119").Value = 34943
$ws.Range("J119").Value = 34943
$ws.Range("L119").Value = 34943
$ws.Range("N119").Value = -44619
# Row 132
$ws.Range("H132").Value = 124750
$ws.Range("J132").Value = 124750
$ws.Range("L132").Value = 124750
$ws.Range("N132").Value = -134870
# Row 134
$ws.Range("H134").Value = 7443.087
$ws.Range("I134").Value = 6216.857
$ws.Range("K134").Value = 18650.571
$ws.Range("M134").Value = -16115.571
# Row 139
$ws.Range("H139").Value = 97999
$ws.Range("J139").Value = 97999
$ws.Range("L139").Value = 97999
$ws.Range("N139").Value = -108279

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2666.6562
$ws.Range("I31").Value = 1530.25
$ws.Range("J31").Value = 3045.4583
$ws.Range("K31").Value = 1530.25
$ws.Range("L31").Value = 3045.4583
$ws.Range("M31").Value = -1235.25
$ws.Range("N31").Value = -3635.4583
# Row 34
$ws.Range("H34").Value = 2666.6562
$ws.Range("I34").Value = 1530.25
$ws.Range("J34").Value = 3045.4583
$ws.Range("K34").Value = 1530.25
$ws.Range("L34").Value = 3045.4583
$ws.Range("M34").Value = -1328.25
$ws.Range("N34").Value = -3449.4583
# Row 99
$ws.Range("H99").Value = 2644.75
$ws.Range("I99").Value = 2499.5
$ws.Range("J99").Value = 2790
$ws.Range("K99").Value = 2499.5
$ws.Range("L99").Value = 2790
$ws.Range("M99").Value = -1001.5
$ws.Range("N99").Value = -5786
# Row 126
$ws.Range("H126").Value = 2644.75
$ws.Range("I126").Value = 2499.5
$ws.Range("J126").Value = 2790
$ws.Range("K126").Value = 7498.5
$ws.Range("L126").Value = 8370
$ws.Range("M126").Value = -5028.5
$ws.Range("N126").Value = -13310

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 494.5
$ws.Range("J2").Value = 120.71429
$ws.Range("L2").Value = 724.28574
$ws.Range("N2").Value = -950.28574
# Row 132
$ws.Range("H132").Value = 7060.067
$ws.Range("I132").Value = 6091.2856
$ws.Range("J132").Value = 7907.75
$ws.Range("K132").Value = 54821.5704
$ws.Range("L132").Value = 71169.75
$ws.Range("M132").Value = -52291.5704
$ws.Range("N132").Value = -76229.75

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 145.90909
$ws.Range("I2").Value = 157.75
$ws.Range("K2").Value = 157.75
$ws.Range("M2").Value = -44.75
# Row 97
$ws.Range("H97").Value = 1801.1515
$ws.Range("I97").Value = 848.7692
$ws.Range("K97").Value = 848.7692
$ws.Range("M97").Value = -352.7692
# Row 135
$ws.Range("H135").Value = 127000
$ws.Range("J135").Value = 127000
$ws.Range("L135").Value = 127000
$ws.Range("N135").Value = -137140

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1287
$ws.Range("J22").Value = 1387.9
$ws.Range("L22").Value = 1387.9
$ws.Range("N22").Value = -1977.9
# Row 27
$ws.Range("H27").Value = 1287
$ws.Range("J27").Value = 1387.9
$ws.Range("L27").Value = 1387.9
$ws.Range("N27").Value = -1601.9
# Row 46
$ws.Range("H46").Value = 1907
$ws.Range("J46").Value = 3099
$ws.Range("L46").Value = 3099
$ws.Range("N46").Value = -3475
# Row 55
$ws.Range("H55").Value = 7144129
$ws.Range("I55").Value = 12500384
$ws.Range("K55").Value = 12500384
$ws.Range("M55").Value = -12500211
# Row 61
$ws.Range("H61").Value = 3493.0715
$ws.Range("I61").Value = 3525.4546
$ws.Range("K61").Value = 3525.4546
$ws.Range("M61").Value = -3323.4546
# Row 100
$ws.Range("H100").Value = 1802
$ws.Range("I100").Value = 1800
$ws.Range("J100").Value = 1804
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 1804
$ws.Range("M100").Value = -1259
$ws.Range("N100").Value = -2886
# Row 113
$ws.Range("H113").Value = 3493.0715
$ws.Range("I113").Value = 3525.4546
$ws.Range("K113").Value = 3525.4546
$ws.Range("M113").Value = -1355.4546
# Row 136
$ws.Range("H136").Value = 2934.8096
$ws.Range("I136").Value = 1526.8
$ws.Range("K136").Value = 4580.4
$ws.Range("M136").Value = -2030.4

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 779023
$ws.Range("I14").Value = 1010930
$ws.Range("J14").Value = 5999.6665
$ws.Range("K14").Value = 1010930
$ws.Range("L14").Value = 5999.6665
$ws.Range("M14").Value = -1010762
$ws.Range("N14").Value = -6335.6665
# Row 100
$ws.Range("H100").Value = 3064.4
$ws.Range("I100").Value = 712.25
$ws.Range("J100").Value = 4632.5
$ws.Range("K100").Value = 1424.5
$ws.Range("L100").Value = 9265
$ws.Range("M100").Value = -883.5
$ws.Range("N100").Value = -10347
# Row 126
$ws.Range("H126").Value = 1154.95
$ws.Range("I126").Value = 1224.4286
$ws.Range("J126").Value = 992.8333
$ws.Range("K126").Value = 3673.2858
$ws.Range("L126").Value = 2978.4999
$ws.Range("M126").Value = -1203.2858
$ws.Range("N126").Value = -7918.4999
# Row 136
$ws.Range("H136").Value = 9937.030000000001
$ws.Range("I136").Value = 8367.223
$ws.Range("J136").Value = 17001.166
$ws.Range("K136").Value = 25101.669
$ws.Range("L136").Value = 51003.49800000001
$ws.Range("M136").Value = -22551.669
$ws.Range("N136").Value = -56103.49800000001

# --- Special case: ARM row 10, M10 cell removed entirely ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M10").ClearContents()
